$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 76; this shifts the existing rows
# 76-87 down to 77-88 (formats/styles carry over automatically).
$ws.Rows(76).Insert()

# Populate the newly inserted row 76 with its data.
$ws.Range("A76").Value = 10
$ws.Range("B76").Value = "Vega Modelo de Temuco"
$ws.Range("C76").Value = "La Araucanía"
$ws.Range("D76").Value = 44637
$ws.Range("E76").Value = 9
$ws.Range("F76").Value = "Fruta"
$ws.Range("G76").Value = 100101
$ws.Range("H76").Value = "Berries"
$ws.Range("I76").Value = 100101001
$ws.Range("J76").Value = "Arándano (blue)"
$ws.Range("K76").Value = "Sin especificar"
$ws.Range("L76").Value = "Primera"
$ws.Range("M76").Value = 80
$ws.Range("N76").Value = 1800
$ws.Range("O76").Value = 1800
$ws.Range("P76").Value = 1800
$ws.Range("Q76").Value = "$/kilo"
$ws.Range("R76").Value = "Región de La Araucanía"
$ws.Range("S76").Value = 1800
$ws.Range("T76").Value = 1
